# Order the cities table by population (descending), per commit message
# "order cities by population". The sheet has a header row (City, Country,
# Population, Area) in row 1 and data in rows 2:13 - everything was already
# sorted by Population descending except Pyeongchang/South Korea, which was
# out of place. Re-sorting the whole table reproduces the intended order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A1:D13")
$sortKey = $ws.Range("C2:C13")

# Key1 = Population column, Order1 = xlDescending (2), Header = xlYes (1)
$dataRange.Sort($sortKey, 2, [System.Type]::Missing, [System.Type]::Missing, 1, [System.Type]::Missing, 1, 1)
